$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Version number content control: "V.1.0" -> "V.1.2"
#    (the plain-text content control that holds the version number,
#    tagged "Version", right after the "Versienummer:" label)
# ------------------------------------------------------------------
foreach ($cc in $d.ContentControls) {
    if ($cc.Tag -eq "Version") {
        $cc.Delete($true)
    }
}

$para = $d.Paragraphs.Item(4)
$para.Range.InsertBefore("V.1.2")

# ------------------------------------------------------------------
# 2. Versiebeheer revision table: append text to the last change note
# ------------------------------------------------------------------
$d.Content.Find.Execute("Wireframes toegevoegd", $true, $false, $false, $false, $false, $true, 1, $false, "Wireframes toegevoegd, navigatiestructuur aangepast", 2)

# ------------------------------------------------------------------
# 3. "Hoofdpagina" heading text -> "Zie pdf"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Hoofdpagina", $true, $false, $false, $false, $false, $true, 1, $false, "Zie pdf", 2)
